$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 154, shifting existing rows 154:175 down to 155:176.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new weekly price entry.
# (Columns A,B,C,E,F,G,H,N,O,Q,R are constant for this product across the sheet.)
$ws.Range("A154").Value = 1
$ws.Range("B154").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C154").Value = "Arica y Parinacota"
$ws.Range("D154").Value = 45142
$ws.Range("E154").Value = 15
$ws.Range("F154").Value = 100112042
$ws.Range("G154").Value = "Locoto"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 160
$ws.Range("K154").Value = 19000
$ws.Range("L154").Value = 20000
$ws.Range("M154").Value = 19500
$ws.Range("N154").Value = "$/caja 20 kilos"
$ws.Range("O154").Value = "Región de Arica y Parinacota"
$ws.Range("P154").Value = 975
$ws.Range("Q154").Value = 20
$ws.Range("R154").Value = "Hortaliza"
